# edit.ps1 - applies the "Added more to motivation section" commit
# to the Word document currently open as $word.ActiveDocument.
#
# Strategy: Word COM (and this iron_native shim) applies Range.InsertXML
# reliably when the target Range spans one or more *complete* paragraphs
# (i.e. includes the trailing paragraph mark). So for every paragraph we
# touch we locate it with Find, grab its full Paragraph.Range (which
# always includes the pilcrow), and replace it wholesale with the
# corrected OOXML for that paragraph (preserving paraId/rsid metadata).

$d = $word.ActiveDocument

function Replace-ParagraphXml {
    param(
        $Para,
        [string]$InnerXml
    )
    $pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
    $pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $full = $pkgHeader + $InnerXml + $pkgFooter
    [void]$Para.Range.InsertXML($full)
}

function Find-Paragraph {
    param([string]$Text)
    $rng = $d.Content
    $ok = $rng.Find.Execute($Text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Could not find text: $Text"
    }
    return $rng.Paragraphs(1)
}

# --- 1. Motivation section: add two new sentences + a following blank paragraph ---
$para = Find-Paragraph "users to access some of the higher-end features without incurring a fee."
$nextPara = $para.Next()
$innerXml = '<w:p w14:paraId="0A7C8011" w14:textId="77777777" w:rsidR="00C33493" w:rsidRPr="00A542C6" w:rsidRDefault="00C33493" w:rsidP="00C33493"><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:tab/><w:t xml:space="preserve">We also believe that this could help enhance the community around campus. Many people play sports, are part of sports clubs, or even just like to stay active in general, and this application could be a way for people to connect with others. </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>An example would be running. If you are not on the track team, there is no real running club on campus, this app could help you find other people around who also like to run. This could end up motivating people to stay active by doing activities together or by making it a fun competition by comparing statistics at the end of each week.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr></w:pPr></w:p>'
Replace-ParagraphXml $nextPara $innerXml

# --- 2. Student Interaction bullet: lastRenderedPageBreak moves here ---
$para = Find-Paragraph "Students will be able to make"
$innerXml = '<w:p w14:paraId="4179C82E" w14:textId="667D6AA1" w:rsidR="00A542C6" w:rsidRDefault="00A56B90" w:rsidP="00C33493"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr><w:spacing w:line="360" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:lastRenderedPageBreak/><w:t>Students will be able to make</w:t></w:r><w:r w:rsidR="00177364"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve"> their own</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00177364"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve">personalized </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve">groups for </w:t></w:r><w:r w:rsidR="00FC746E"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve">both </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve">campus clubs and friends to </w:t></w:r><w:r w:rsidR="00763E15"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>track activities</w:t></w:r><w:r w:rsidR="00F31E19"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>, allowing others to participate</w:t></w:r><w:r w:rsidR="00FC746E"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>.</w:t></w:r><w:r w:rsidR="00AD146C"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve"> The Vulcan Activity Tracker will display </w:t></w:r><w:r w:rsidR="00F31E19"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>what</w:t></w:r><w:r w:rsidR="00FC746E"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00F31E19"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>activities were done as a group</w:t></w:r><w:r w:rsidR="00C33493"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve"> and the members involved</w:t></w:r><w:r w:rsidR="00C74BE4"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve">. </w:t></w:r></w:p>'
Replace-ParagraphXml $para $innerXml

# --- 3. Implementation Techniques heading: lastRenderedPageBreak removed ---
$para = Find-Paragraph "Implementation "
$innerXml = '<w:p w14:paraId="4D8864C0" w14:textId="10118C36" w:rsidR="00C8378F" w:rsidRPr="00A542C6" w:rsidRDefault="00C8378F" w:rsidP="00C33493"><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:i/><w:iCs/></w:rPr></w:pPr><w:r w:rsidRPr="00A542C6"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:i/><w:iCs/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t xml:space="preserve">Implementation </w:t></w:r><w:r w:rsidR="00801366" w:rsidRPr="00A542C6"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:i/><w:iCs/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>T</w:t></w:r><w:r w:rsidRPr="00A542C6"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:i/><w:iCs/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>echniques</w:t></w:r></w:p>'
Replace-ParagraphXml $para $innerXml

# --- 4. VS Code heading: lastRenderedPageBreak moves here ---
$para = Find-Paragraph "VS Code"
$innerXml = '<w:p w14:paraId="50E62438" w14:textId="7D60D85E" w:rsidR="00281D61" w:rsidRDefault="00281D61" w:rsidP="008F5A06"><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:ind w:firstLine="720"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:i/><w:iCs/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:i/><w:iCs/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:lastRenderedPageBreak/><w:t>VS Code</w:t></w:r></w:p>'
Replace-ParagraphXml $para $innerXml

# --- 5. Potential Users heading: lastRenderedPageBreak removed ---
$para = Find-Paragraph "Potential Users"
$innerXml = '<w:p w14:paraId="63B6BF27" w14:textId="74CD33E9" w:rsidR="00C8378F" w:rsidRPr="00A542C6" w:rsidRDefault="00C8378F" w:rsidP="00C33493"><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:i/><w:iCs/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r w:rsidRPr="00A542C6"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:i/><w:iCs/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>Potential Users</w:t></w:r></w:p>'
Replace-ParagraphXml $para $innerXml

# --- 6. Routes paragraph: split run, add lastRenderedPageBreak on 2nd half ---
$para = Find-Paragraph "rotary park loop"
$innerXml = '<w:p w14:paraId="4FFE0A27" w14:textId="5FFDA62E" w:rsidR="00C8378F" w:rsidRPr="00B169D7" w:rsidRDefault="00C8378F" w:rsidP="00C33493"><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">Routes </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve">– the program will have a feature that will allow a user to select the route they took if their activity was outdoors. There will be a list of available routes along with their mileage. Users can select how many times they did the route, or if they did a </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">second route within the activity (an example would be going into the rotary park loop and then coming back out). </w:t></w:r></w:p>'
Replace-ParagraphXml $para $innerXml

# --- 7. HTML frontend paragraph: merge ' reduce ' / 'all of' / ' the code...' runs, drop proofErr ---
$para = Find-Paragraph "all of"
$innerXml = '<w:p w14:paraId="6FF404C9" w14:textId="4C754F3B" w:rsidR="00C8378F" w:rsidRPr="00155D8C" w:rsidRDefault="00C8378F" w:rsidP="00C33493"><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:i/><w:iCs/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve">HTML and CSS will be used for our frontend framework of this web-based application. We will make use of base templates </w:t></w:r><w:r w:rsidR="00C33493"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>to</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve"> reduce all of the code to one place and to make routing to another page easier. CSS will be used for the style of the website</w:t></w:r><w:r w:rsidR="00884577"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t xml:space="preserve">. </w:t></w:r></w:p>'
Replace-ParagraphXml $para $innerXml

Write-Output "done"
